$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Location / FPK12 Exam / FPK12 Schedule values (parallel-execution test data refresh)
$ws.Range("A2").Value = "ECLocation253798"
$ws.Range("H2").Value = "FPK12Exam71442"
$ws.Range("I2").Value = "FPK12Schedule410184"

# Rows 3-5 - Tenant ID column values
$ws.Range("E3").Value = "917893"
$ws.Range("E4").Value = "604713"
$ws.Range("E5").Value = "547406"
